# Fourth commit on 02-02-2026
# Remove the "2" page-number textbox and the "Running First java program"
# title textbox from slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Delete from the highest shape index down so indices of the shapes we
# still need to remove don't shift underneath us.
$s.Shapes.Item("TextBox 5").Delete()
$s.Shapes.Item("TextBox 2").Delete()
